$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column M (the "alcohol" data column that held running totals);
# the following column (old N) shifts left to become the new column M.
$ws.Range("M1").EntireColumn.Delete() | Out-Null

# Move/restore the active selection onto the (now last) column M, row 1.
$ws.Range("M1").Select() | Out-Null
